$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Make a payment" query row (row 22)
$ws.Range("A22").Value = "Payment"
$ws.Range("B22").Value = "Make a payment- Once off - Enter amount"
$ws.Range("C22").Value = "Preference.yaml"
$ws.Range("D22").Value = "GET"
$ws.Range("E22").Value = "limits"
$ws.Range("F22").Value = "1. How to calculate instant payment charges for specific payment.`n2. How the daily payment limit will be set. (Whether it is selected account level or user profile dependent)"
$ws.Range("F22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 46.8

# Leave a review comment on the "limits" endpoint cell, like the author did
$comment = $ws.Range("E22").AddComment("Pawan Jashnani:`ni")

# Select the new row, as the author would have left the sheet positioned there
$ws.Range("A22").Select() | Out-Null
